$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2023 column (K) to the table, matching the look/format of
# the preceding 2022 column (J), then fill in the actual figures.

$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)

$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)

$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 904.5
$ws.Range("K5").Value = 662.6
$ws.Range("K6").Value = 1147.2
